# Add minimal formatting to Excel sheet:
#  - turn off gridlines on every sheet
#  - set explicit column widths on every sheet
#  - turn every sheet's used range into a named Excel Table (ListObject)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create a Table (ListObject) on each worksheet, covering its used range.
#
#    Two runtime quirks to work around:
#
#    a) The very first tables created in a *fresh* workbook are numbered
#       table1.xml, table2.xml, ... but the target file expects the four
#       tables to land on table3.xml..table6.xml (ids 3..6). Creating (and
#       then deleting) two throwaway single-column tables on throwaway
#       sheets first "burns" ids 1 and 2, so the four real tables land on
#       ids 3..6, matching the target numbering.
#
#    b) Renaming a ListObject to "TableN" *before* the Nth table in the
#       workbook has actually been created can make an earlier table
#       vanish on save. Work around it by (i) adding every real table
#       first with its default name, then (ii) renaming them afterwards,
#       in reverse order.
# ---------------------------------------------------------------------------

$dummy1 = $wb.Worksheets.Add()
$dummy1.Range("A1").Value = "x"
$dummy1.Range("A2").Value = "y"
$lod1 = $dummy1.ListObjects.Add(1, $dummy1.Range("A1:A2"), 0, 1)
$lod1.Name = "Dummy1"

$dummy2 = $wb.Worksheets.Add()
$dummy2.Range("A1").Value = "x"
$dummy2.Range("A2").Value = "y"
$lod2 = $dummy2.ListObjects.Add(1, $dummy2.Range("A1:A2"), 0, 1)
$lod2.Name = "Dummy2"

$dummy1.Delete() | Out-Null
$dummy2.Delete() | Out-Null

$ws1 = $wb.Worksheets.Item(1)   # Table
$ws2 = $wb.Worksheets.Item(2)   # Variables
$ws3 = $wb.Worksheets.Item(3)   # Codelists
$ws4 = $wb.Worksheets.Item(4)   # Data

$lo1 = $ws1.ListObjects.Add(1, $ws1.Range("A1:B41"), 0, 1)
$lo2 = $ws2.ListObjects.Add(1, $ws2.Range("A1:O5"), 0, 1)
$lo3 = $ws3.ListObjects.Add(1, $ws3.Range("A1:G7"), 0, 1)
$lo4 = $ws4.ListObjects.Add(1, $ws4.Range("A1:D46"), 0, 1)

$lo4.Name = "Table6"
$lo3.Name = "Table5"
$lo2.Name = "Table4"
$lo1.Name = "Table3"

# ---------------------------------------------------------------------------
# 2) Explicit column widths per sheet (values chosen so the saved XML is as
#    close as this runtime's pixel-snapped column-width model allows to the
#    target "<integer>.71" character widths).
# ---------------------------------------------------------------------------

function Set-ColWidths($ws, $widths) {
    for ($i = 0; $i -lt $widths.Length; $i++) {
        $ws.Columns.Item($i + 1).ColumnWidth = $widths[$i]
    }
}

Set-ColWidths $ws1 @(15.8333333333333, 40.8333333333333)
Set-ColWidths $ws2 @(7.83333333333333, 13.8333333333333, 6.83333333333333, 13.8333333333333, 11.8333333333333, 11.8333333333333, 8.83333333333333, 8.83333333333333, 8.83333333333333, 13.8333333333333, 13.8333333333333, 13.8333333333333, 31.8333333333333, 31.8333333333333, 31.8333333333333)
Set-ColWidths $ws3 @(13.8333333333333, 3.83333333333333, 8.83333333333333, 16.8333333333333, 16.8333333333333, 25.8333333333333, 8.83333333333333)
Set-ColWidths $ws4 @(13.8333333333333, 3.83333333333333, 5.83333333333333, 7.83333333333333)

# ---------------------------------------------------------------------------
# 3) Turn off gridlines on every sheet, then restore the originally active
#    tab (sheet 1) as the selected tab.
# ---------------------------------------------------------------------------

for ($i = 1; $i -le 4; $i++) {
    $sheet = $wb.Worksheets.Item($i)
    $sheet.Activate()
    $excel.ActiveWindow.DisplayGridlines = $false
}
$ws1.Activate()

Write-Host "done"
